$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G (header "K") values regenerated from std/mean calc (s_vals)
$kValues = @{
    2 = 2
    3 = 1
    4 = 0
    5 = 1
    6 = 0
    7 = 2
    8 = 2
    9 = 0
    10 = 0
    11 = 3
    12 = 1
    13 = 0
    14 = 0
    15 = 0
    16 = 1
    17 = 1
    18 = 2
    19 = 0
    20 = 1
    21 = 2
    22 = 1
    23 = 1
    24 = 1
    25 = 2
    26 = 1
    27 = 1
    28 = 1
    30 = 1
    31 = 1
    32 = 1
    33 = 1
    34 = 1
    35 = 2
    36 = 1
    37 = 0
    38 = 0
    39 = 1
    42 = 0
    43 = 0
    44 = 3
    45 = 0
    46 = 1
    47 = 0
    48 = 0
    49 = 1
    50 = 1
    51 = 2
    52 = 2
    53 = 1
    54 = 1
    55 = 0
    56 = 0
    57 = 0
    58 = 0
    59 = 1
    60 = 3
    61 = 0
    62 = 1
    63 = 1
    64 = 2
    65 = 0
    66 = 0
    67 = 1
    68 = 3
    69 = 1
    70 = 1
    71 = 2
    72 = 1
    73 = 1
    74 = 1
    75 = 1
    76 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}

